# Updates the cryptos list worksheet with latest price/volume snapshot values.
# For the "Price" (D) column we temporarily force a Text number format before
# assigning the value so Excel keeps values such as "90.194.55" or "0.365"
# as literal text instead of silently re-interpreting them as numbers/dates.
# We then restore the "Normal" style so no stray number-format/style index is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.194.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.100.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -12.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.365"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.098.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.720"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.198"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.882.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.624.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("B18").Value = "SuiNetwork"
$ws.Range("C18").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.094.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000214"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "434.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.273.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.193"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.152"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "498.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +56.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0884"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -6.13%  "
$ws.Range("E46").Value = "  -6.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.681"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "149.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.11%  "

